{"js": "// Insert the text \"Angular Projects\" into the first (empty) paragraph of the\n// document body, matching the diff which turns the leading `<w:p/>` into a\n// paragraph containing a single run with that text.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.insertText(\"Angular Projects\", Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Insert the text \"Angular Projects\" into the first (empty) paragraph of the\n# document body, matching the diff which turns the leading empty paragraph\n# into a paragraph containing the text \"Angular Projects\".\n\n$d = $word.ActiveDocument\n$firstPara = $d.Paragraphs.Item(1)\n$firstPara.Range.InsertBefore(\"Angular Projects\")\n"}
